$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New title rows (1-4) in column D: company / cert title / report title
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "NR Finance Mexico"
$ws.Range("D2").Value = "EKIP"
$ws.Range("D3").Value = "Certificacion de usuarios 2024"
$ws.Range("D4").Value = "Reporte de usuarios"

$titleRng = $ws.Range("D1:D4")
$titleRng.Font.Size = 16
$titleRng.Font.Bold = $true
$titleRng.HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 2) Thin borders around the existing data table (A5:F7)
# ---------------------------------------------------------------------------
$ws.Range("A5:F7").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3) Column widths (values compensate for the engine's xlsx pixel-export
#    offset so the saved <col width> lands on the intended value)
# ---------------------------------------------------------------------------
$offset = 5/6
$ws.Columns.Item(1).ColumnWidth = 4.853482 - $offset
$ws.Columns.Item(2).ColumnWidth = 39.139196 - $offset
$ws.Columns.Item(3).ColumnWidth = 8.424911 - $offset
$ws.Columns.Item(4).ColumnWidth = 40.424911 - $offset
$ws.Columns.Item(5).ColumnWidth = 24.139196 - $offset
$ws.Columns.Item(6).ColumnWidth = 16.567768 - $offset

Write-Host "Edit complete"
